$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.601795673370361
$ws.Range("B1").Value = 1.719259977340698
$ws.Range("C1").Value = 1.967326045036316
$ws.Range("D1").Value = 3.174352645874023
$ws.Range("E1").Value = 3.781453132629395
